$wb = $excel.ActiveWorkbook

# --- Content fix: sheet "机柜与插接箱" template placeholders ---
# {.boxIndexa} -> {.outletIda}  (cell J3)
# {.boxIndexb} -> {.outletIdb}  (cell M3)
$wsBox = $wb.Worksheets.Item("机柜与插接箱")
$wsBox.Range("J3").Value = "{.outletIda}"
$wsBox.Range("M3").Value = "{.outletIdb}"

# --- View-state fix: restore/update the saved selection on each sheet ---
# Sheet "机柜与插接箱": selection moves from I33 to O5
$wsBox.Activate()
$wsBox.Range("O5").Select()

# Sheet "机柜与PDU": selection moves from O11 to H9, and this sheet stays
# the active/tab-selected sheet (matches the original workbook state).
$wsPdu = $wb.Worksheets.Item("机柜与PDU")
$wsPdu.Activate()
$wsPdu.Range("H9").Select()
